$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("G2").Value = 73
$ws.Range("H2").Value = 92
$ws.Range("D3").Value = 109
$ws.Range("G3").Value = 113
$ws.Range("K3").Value = 183
$ws.Range("F4").Value = 5
$ws.Range("C6").Value = 392
$ws.Range("E6").Value = 357
$ws.Range("F6").Value = 428
$ws.Range("G6").Value = 388
$ws.Range("H6").Value = 373
$ws.Range("J6").Value = 335
$ws.Range("K6").Value = 422
$ws.Range("C7").Value = 524
$ws.Range("D7").Value = 525
$ws.Range("E7").Value = 537
$ws.Range("F7").Value = 605
$ws.Range("G7").Value = 578
$ws.Range("H7").Value = 594
$ws.Range("J7").Value = 628
$ws.Range("K7").Value = 745

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("C6").Value = 29
$ws.Range("F6").Value = 47
$ws.Range("C7").Value = 32
$ws.Range("F7").Value = 55

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("F4").Value = 2
$ws.Range("C6").Value = 28
$ws.Range("K6").Value = 35
$ws.Range("C7").Value = 33
$ws.Range("F7").Value = 40
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 9

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("G2").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("D6").Value = 5
$ws.Range("G6").Value = 8

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("F6").Value = 32
$ws.Range("F7").Value = 43

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 8

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("G2").Value = 4
$ws.Range("C19").Value = 9
$ws.Range("K19").Value = 30
$ws.Range("F28").Value = 43
$ws.Range("K29").Value = 17
$ws.Range("J30").Value = 5
$ws.Range("C32").Value = 32
$ws.Range("F32").Value = 55
$ws.Range("C36").Value = 33
$ws.Range("F36").Value = 40
$ws.Range("K36").Value = 54
$ws.Range("D41").Value = 5
$ws.Range("G41").Value = 8
$ws.Range("J42").Value = 8
$ws.Range("K53").Value = 80
$ws.Range("G62").Value = 6
$ws.Range("H63").Value = 2
$ws.Range("K75").Value = 3
$ws.Range("H76").Value = 16
$ws.Range("E82").Value = 8
$ws.Range("G88").Value = 9
$ws.Range("J91").Value = 4
$ws.Range("C98").Value = 524
$ws.Range("D98").Value = 525
$ws.Range("E98").Value = 537
$ws.Range("F98").Value = 605
$ws.Range("G98").Value = 578
$ws.Range("H98").Value = 594
$ws.Range("J98").Value = 628
$ws.Range("K98").Value = 745

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 4

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("H5").Value = 12
$ws.Range("H6").Value = 16

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("G3").Value = 2
$ws.Range("G7").Value = 6

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K5").Value = 13
$ws.Range("K6").Value = 17

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("H4").Value = 7
$ws.Range("H5").Value = 8

$ws = $wb.Worksheets.Item("New City")
$ws.Range("G2").Value = 2
$ws.Range("G5").Value = 2

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("C6").Value = 8
$ws.Range("K6").Value = 15
$ws.Range("C7").Value = 9
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("F3").Value = 2
$ws.Range("F6").Value = 4

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I5").Value = 3
$ws.Range("I6").Value = 5

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J4").Value = 3
$ws.Range("J5").Value = 3

